$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.187.38'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.642.93'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.36'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.95'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.140'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.20%  '
$ws.Range("E10").Value = '  -1.40%  '
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.97'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.123.52'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.103.76'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.626.16'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.51%  '
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '363.32'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.42'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.02%  '
$ws.Range("E22").Value = '  -2.52%  '
$ws.Range("E23").Value = '  -3.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.17'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.05%  '
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("E26").Value = '  -3.62%  '
$ws.Range("E27").Value = '  +7.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.776.70'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000104'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '555.02'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.61%  '
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.85'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  -2.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.54'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.57'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.15%  '
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("E40").Value = '  -3.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.32'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("E42").Value = '  +3.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.85'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("E44").Value = '  -2.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '158.71'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.05'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("E49").Value = '  -2.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0781'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("E51").Value = '  -0.93%  '
